$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("All Trades", "leadlag")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A6").Value = 5
    $ws.Range("B6").Formula = '="2026-02-16"'
    $ws.Range("C6").Value = "21:51:21"
    $ws.Range("D6").Value = "leadlag"
    $ws.Range("E6").Value = "DOWN"
    $ws.Range("F6").Value = 68306.38
    $ws.Range("G6").Formula = '=""'
    $ws.Range("H6").Value = "OPEN"
    $ws.Range("I6").Value = 0
    $ws.Range("J6").Value = 0
    $ws.Range("K6").Value = 100
    $ws.Range("L6").Value = 0.75
    $ws.Range("M6").Value = "Coinbase leading with -0.117% move"
    $ws.Range("N6").Formula = '=""'
    $ws.Range("O6").Value = 0
}
